# Apply the "Method column into Type" update to biovolume_method.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Bryozoa / Polyp): remove the bryozoan colonial-specimen note from column F
$ws.Range("F7").ClearContents()

# Row 31 (Ophiuroidea): correct tense of the arms/discs note in column F
$ws.Range("F31").Value = "Both arms and discs are assumed to be cylindrical"

# Row 36 (Polychaeta): the Type is now "Cylinder" instead of "LWR", and the
# associated LWR coefficient value no longer applies
$ws.Range("C36").Value = "Cylinder"
$ws.Range("D36").ClearContents()

# Reflect final selection / scroll position left by the editor
$ws.Range("A36:XFD36").Select()
$excel.ActiveWindow.ScrollRow = 17

$wb.Save()
